$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-29"

# Update the September row label text
$ws.Range("A10").Value = "September (through 09-29)"

# Row 9 (August) - only 2021/2022 columns updated
$ws.Range("H9").Value = 159
$ws.Range("I9").Value = 165

# Row 10 (September) - all year columns updated
$ws.Range("B10").Value = 30
$ws.Range("C10").Value = 44
$ws.Range("D10").Value = 74
$ws.Range("E10").Value = 54
$ws.Range("F10").Value = 71
$ws.Range("G10").Value = 112
$ws.Range("H10").Value = 173
$ws.Range("I10").Value = 139

# Row 11 (Total) - all year columns updated
$ws.Range("B11").Value = 224
$ws.Range("C11").Value = 425
$ws.Range("D11").Value = 625
$ws.Range("E11").Value = 544
$ws.Range("F11").Value = 420
$ws.Range("G11").Value = 896
$ws.Range("H11").Value = 1242
$ws.Range("I11").Value = 1273
